$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 626.2727
$ws.Range("I6").Value = 626.2727
$ws.Range("K6").Value = 1878.8181
$ws.Range("M6").Value = -1766.8181
$ws.Range("H29").Value = 978.0
$ws.Range("I29").Value = 978.0
$ws.Range("K29").Value = 2934.0
$ws.Range("M29").Value = -2653.0
$ws.Range("H38").Value = 5695.6
$ws.Range("I38").Value = 344.0
$ws.Range("J38").Value = 9263.333
$ws.Range("K38").Value = 1032.0
$ws.Range("L38").Value = 27789.999
$ws.Range("M38").Value = -660.0
$ws.Range("N38").Value = -28533.999
$ws.Range("H58").Value = 4565.875
$ws.Range("I58").Value = 424.0
$ws.Range("J58").Value = 8070.5386
$ws.Range("K58").Value = 1272.0
$ws.Range("L58").Value = 24211.6158
$ws.Range("M58").Value = -1122.0
$ws.Range("N58").Value = -24511.6158
$ws.Range("H87").Value = 20875.334
$ws.Range("J87").Value = 20875.334
$ws.Range("L87").Value = 20875.334
$ws.Range("N87").Value = -23371.334
$ws.Range("H90").Value = 20875.334
$ws.Range("J90").Value = 20875.334
$ws.Range("L90").Value = 62626.00199999999
$ws.Range("N90").Value = -75106.002
$ws.Range("H100").Value = 33337496.0
$ws.Range("I100").Value = 66669330.0
$ws.Range("K100").Value = 66669330.0
$ws.Range("M100").Value = -66668789.0
$ws.Range("H112").Value = 1350.5
$ws.Range("J112").Value = 1364.6154
$ws.Range("L112").Value = 4093.8462
$ws.Range("N112").Value = -6309.8462
$ws.Range("H125").Value = 754.0
$ws.Range("I125").Value = 600.3333
$ws.Range("J125").Value = 885.7143
$ws.Range("K125").Value = 5402.9997
$ws.Range("L125").Value = 7971.428699999999
$ws.Range("M125").Value = -2942.9997
$ws.Range("N125").Value = -12891.4287
$ws.Range("H132").Value = 21335240.0
$ws.Range("I132").Value = 23582950.0
$ws.Range("J132").Value = 2004932.4
$ws.Range("K132").Value = 70748850.0
$ws.Range("L132").Value = 6014797.199999999
$ws.Range("M132").Value = -70746320.0
$ws.Range("N132").Value = -6019857.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3660.0278
$ws.Range("I32").Value = 3495.4067
$ws.Range("J32").Value = 4407.154
$ws.Range("K32").Value = 3495.4067
$ws.Range("L32").Value = 4407.154
$ws.Range("M32").Value = -3208.4067
$ws.Range("N32").Value = -4981.154
$ws.Range("H81").Value = 33500.0
$ws.Range("J81").Value = 33500.0
$ws.Range("L81").Value = 33500.0
$ws.Range("N81").Value = -35496.0
$ws.Range("H84").Value = 33500.0
$ws.Range("J84").Value = 33500.0
$ws.Range("L84").Value = 100500.0
$ws.Range("N84").Value = -110484.0
$ws.Range("H97").Value = 1555.5555
$ws.Range("I97").Value = 1515.0
$ws.Range("J97").Value = 1636.6666
$ws.Range("K97").Value = 1515.0
$ws.Range("L97").Value = 1636.6666
$ws.Range("M97").Value = -1019.0
$ws.Range("N97").Value = -2628.6666
$ws.Range("H137").Value = 38786.0
$ws.Range("J137").Value = 38786.0
$ws.Range("L137").Value = 38786.0
$ws.Range("N137").Value = -48986.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1202.48
$ws.Range("I107").Value = 1102.125
$ws.Range("J107").Value = 1380.8889
$ws.Range("K107").Value = 1102.125
$ws.Range("L107").Value = 1380.8889
$ws.Range("M107").Value = 817.875
$ws.Range("N107").Value = -5220.8889
$ws.Range("H112").Value = 29992.5
$ws.Range("J112").Value = 29992.5
$ws.Range("L112").Value = 29992.5
$ws.Range("N112").Value = -32946.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 235413.27
$ws.Range("I31").Value = 966341.0
$ws.Range("J31").Value = 2845.3408
$ws.Range("K31").Value = 966341.0
$ws.Range("L31").Value = 2845.3408
$ws.Range("M31").Value = -966046.0
$ws.Range("N31").Value = -3435.3408
$ws.Range("H34").Value = 235413.27
$ws.Range("I34").Value = 966341.0
$ws.Range("J34").Value = 2845.3408
$ws.Range("K34").Value = 966341.0
$ws.Range("L34").Value = 2845.3408
$ws.Range("M34").Value = -966139.0
$ws.Range("N34").Value = -3249.3408
$ws.Range("H99").Value = 4292.2354
$ws.Range("I99").Value = 2449.25
$ws.Range("J99").Value = 5930.4443
$ws.Range("K99").Value = 2449.25
$ws.Range("L99").Value = 5930.4443
$ws.Range("M99").Value = -951.25
$ws.Range("N99").Value = -8926.4443
$ws.Range("H105").Value = 1264.1111
$ws.Range("I105").Value = 1100.9546
$ws.Range("J105").Value = 1982.0
$ws.Range("K105").Value = 1100.9546
$ws.Range("L105").Value = 1982.0
$ws.Range("M105").Value = 646.0454
$ws.Range("N105").Value = -5476.0
$ws.Range("H115").Value = 24833.0
$ws.Range("I115").Value = 24799.0
$ws.Range("J115").Value = 24850.0
$ws.Range("K115").Value = 24799.0
$ws.Range("L115").Value = 24850.0
$ws.Range("M115").Value = -23624.0
$ws.Range("N115").Value = -27200.0
$ws.Range("H126").Value = 4292.2354
$ws.Range("I126").Value = 2449.25
$ws.Range("J126").Value = 5930.4443
$ws.Range("K126").Value = 7347.75
$ws.Range("L126").Value = 17791.3329
$ws.Range("M126").Value = -4877.75
$ws.Range("N126").Value = -22731.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2312.75
$ws.Range("I17").Value = 501.0
$ws.Range("K17").Value = 1503.0
$ws.Range("M17").Value = -1334.0
$ws.Range("H34").Value = 12355.158
$ws.Range("J34").Value = 9681.5
$ws.Range("L34").Value = 29044.5
$ws.Range("N34").Value = -29212.5
$ws.Range("H68").Value = 4525.857
$ws.Range("I68").Value = 1625.25
$ws.Range("J68").Value = 5385.2964
$ws.Range("K68").Value = 4875.75
$ws.Range("L68").Value = 16155.8892
$ws.Range("M68").Value = -4064.75
$ws.Range("N68").Value = -17777.8892
$ws.Range("H71").Value = 4525.857
$ws.Range("I71").Value = 1625.25
$ws.Range("J71").Value = 5385.2964
$ws.Range("K71").Value = 14627.25
$ws.Range("L71").Value = 48467.6676
$ws.Range("M71").Value = -10571.25
$ws.Range("N71").Value = -56579.6676
$ws.Range("H112").Value = 701299.6
$ws.Range("I112").Value = 499.33334
$ws.Range("J112").Value = 1752500.0
$ws.Range("K112").Value = 1498.00002
$ws.Range("L112").Value = 5257500.0
$ws.Range("M112").Value = -390.0000199999999
$ws.Range("N112").Value = -5259716.0
$ws.Range("H123").Value = 4307.6665
$ws.Range("I123").Value = 4000.0
$ws.Range("J123").Value = 4923.0
$ws.Range("K123").Value = 12000.0
$ws.Range("L123").Value = 14769.0
$ws.Range("M123").Value = -9550.0
$ws.Range("N123").Value = -19669.0
$ws.Range("H131").Value = 836.45
$ws.Range("I131").Value = 520.5
$ws.Range("J131").Value = 856.617
$ws.Range("K131").Value = 1561.5
$ws.Range("L131").Value = 2569.851
$ws.Range("M131").Value = 3478.5
$ws.Range("N131").Value = -12649.851
$ws.Range("H133").Value = 9899.667
$ws.Range("I133").Value = 14999.333
$ws.Range("K133").Value = 44997.999
$ws.Range("M133").Value = -39937.999
$ws.Range("H134").Value = 2641.5862
$ws.Range("I134").Value = 1630.3
$ws.Range("J134").Value = 4888.8887
$ws.Range("K134").Value = 4890.9
$ws.Range("L134").Value = 14666.6661
$ws.Range("M134").Value = 179.1000000000004
$ws.Range("N134").Value = -24806.6661
$ws.Range("H139").Value = 914.0
$ws.Range("I139").Value = 751.5789
$ws.Range("K139").Value = 2254.7367
$ws.Range("M139").Value = 2885.2633

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1889.6666
$ws.Range("I113").Value = 2080.1667
$ws.Range("J113").Value = 1508.6666
$ws.Range("K113").Value = 2080.1667
$ws.Range("L113").Value = 1508.6666
$ws.Range("M113").Value = 89.83329999999978
$ws.Range("N113").Value = -5848.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3628.9583
$ws.Range("I7").Value = 1987.2222
$ws.Range("J7").Value = 4614.0
$ws.Range("K7").Value = 1987.2222
$ws.Range("L7").Value = 4614.0
$ws.Range("M7").Value = -1875.2222
$ws.Range("N7").Value = -4838.0
$ws.Range("H55").Value = 119.125
$ws.Range("I55").Value = 100.2
$ws.Range("J55").Value = 150.66667
$ws.Range("K55").Value = 100.2
$ws.Range("L55").Value = 150.66667
$ws.Range("M55").Value = 72.8
$ws.Range("N55").Value = -496.66667
$ws.Range("H126").Value = 3628.9583
$ws.Range("I126").Value = 1987.2222
$ws.Range("J126").Value = 4614.0
$ws.Range("K126").Value = 5961.6666
$ws.Range("L126").Value = 13842.0
$ws.Range("M126").Value = -3491.6666
$ws.Range("N126").Value = -18782.0
$ws.Range("H132").Value = 3609.8538
$ws.Range("I132").Value = 2616.2083
$ws.Range("J132").Value = 5012.647
$ws.Range("K132").Value = 7848.624899999999
$ws.Range("L132").Value = 15037.941
$ws.Range("M132").Value = -5318.624899999999
$ws.Range("N132").Value = -20097.941
$ws.Range("H136").Value = 3065.0454
$ws.Range("I136").Value = 1223.2084
$ws.Range("J136").Value = 5275.25
$ws.Range("K136").Value = 3669.6252
$ws.Range("L136").Value = 15825.75
$ws.Range("M136").Value = -1119.6252
$ws.Range("N136").Value = -20925.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 549.5
$ws.Range("I100").Value = 198.0
$ws.Range("J100").Value = 666.6667
$ws.Range("K100").Value = 396.0
$ws.Range("L100").Value = 1333.3334
$ws.Range("M100").Value = 145.0
$ws.Range("N100").Value = -2415.3334
$ws.Range("H125").Value = 37361.0
$ws.Range("J125").Value = 37361.0
$ws.Range("L125").Value = 37361.0
$ws.Range("N125").Value = -47201.0
$ws.Range("H132").Value = 1991.6786
$ws.Range("I132").Value = 857.6818
$ws.Range("K132").Value = 2573.0454
$ws.Range("M132").Value = -43.04539999999997
$ws.Range("H137").Value = 43105.715
$ws.Range("J137").Value = 43105.715
$ws.Range("L137").Value = 43105.715
$ws.Range("N137").Value = -53305.715
$ws.Range("H139").Value = 46364.9
$ws.Range("J139").Value = 46364.9
$ws.Range("L139").Value = 46364.9
$ws.Range("N139").Value = -56644.9
